$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.347730641930184
$ws.Range("C2").Value = 0.06767891441089091
$ws.Range("E2").Value = 0.4143064095259632
$ws.Range("F2").Value = 0.4443680307746121
$ws.Range("G2").Value = 0.6351006585564392
$ws.Range("H2").Value = 0.7402179449673838
$ws.Range("I2").Value = 0.5824378794838729
$ws.Range("K2").Value = 0.3881189099243443
$ws.Range("B3").Value = 0.3074721832826981
$ws.Range("C3").Value = 0.05895372407800892
$ws.Range("E3").Value = 0.3614587690521915
$ws.Range("F3").Value = 0.3878228170618172
$ws.Range("G3").Value = 0.6350068731898375
$ws.Range("H3").Value = 0.7455245927634166
$ws.Range("I3").Value = 0.5879938798843796
$ws.Range("K3").Value = 0.3408481337984881
$ws.Range("B4").Value = 0.2827769100592263
$ws.Range("C4").Value = 0.05358639945231403
$ws.Range("E4").Value = 0.3291050374318303
$ws.Range("F4").Value = 0.3531389305169483
$ws.Range("G4").Value = 0.6355336757560366
$ws.Range("H4").Value = 0.7492282096347935
$ws.Range("I4").Value = 0.5918251646278136
$ws.Range("K4").Value = 0.311831940487906
$ws.Range("B5").Value = 0.2727195437841772
$ws.Range("C5").Value = 0.05139655236725105
$ws.Range("E5").Value = 0.3159426302232475
$ws.Range("F5").Value = 0.3390132514313251
$ws.Range("G5").Value = 0.6358944618867213
$ws.Range("H5").Value = 0.7508491611921642
$ws.Range("I5").Value = 0.5934916489649567
$ws.Range("K5").Value = 0.3000097636492853
$ws.Range("B6").Value = 0.2710499046047801
$ws.Range("C6").Value = 0.05103276892708664
$ws.Range("E6").Value = 0.3137582969655313
$ws.Range("F6").Value = 0.336668177824194
$ws.Range("G6").Value = 0.6359631702979414
$ws.Range("H6").Value = 0.7511250570522776
$ws.Range("I6").Value = 0.5937747105028528
$ws.Range("K6").Value = 0.2980468331535349
$ws.Range("B7").Value = 0.282641247658546
$ws.Range("C7").Value = 0.05355687714019552
$ws.Range("E7").Value = 0.3289274379711742
$ws.Range("F7").Value = 0.3529483938344953
$ws.Range("G7").Value = 0.6355379509673185
$ws.Range("H7").Value = 0.7492496184723336
$ws.Range("I7").Value = 0.5918472139530166
$ws.Range("K7").Value = 0.3116724936242008
$ws.Range("B8").Value = 0.3338447527399637
$ws.Range("C8").Value = 0.0646724866894175
$ws.Range("E8").Value = 0.3960636951937886
$ws.Range("F8").Value = 0.4248636149813478
$ws.Range("G8").Value = 0.6349465244178134
$ws.Range("H8").Value = 0.7419550828865908
$ws.Range("I8").Value = 0.5842662155477676
$ws.Range("K8").Value = 0.3718182651463735
$ws.Range("B9").Value = 0.434437214403431
$ws.Range("C9").Value = 0.08639617459417082
$ws.Range("E9").Value = 0.5285616984483426
$ws.Range("F9").Value = 0.5661985755041457
$ws.Range("G9").Value = 0.6384623835381547
$ws.Range("H9").Value = 0.7311968415709771
$ws.Range("I9").Value = 0.5727485764559717
$ws.Range("K9").Value = 0.4898322340060872
$ws.Range("B10").Value = 0.5084556700172698
$ws.Range("C10").Value = 0.1023207025986608
$ws.Range("E10").Value = 0.6265573754423315
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("G10").Value = 0.6439509584739369
$ws.Range("H10").Value = 0.725473101226811
$ws.Range("I10").Value = 0.5663519106988133
$ws.Range("K10").Value = 0.5765919757120628
$ws.Range("B11").Value = 0.5421541864637902
$ws.Range("C11").Value = 0.1095595078413965
$ws.Range("E11").Value = 0.6713090427410719
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("G11").Value = 0.647090740721822
$ws.Range("H11").Value = 0.7233467008295378
$ws.Range("I11").Value = 0.563895562598077
$ws.Range("K11").Value = 0.6160767047601325
$ws.Range("B12").Value = 0.5549187638534079
$ws.Range("C12").Value = 0.1123000273538253
$ws.Range("E12").Value = 0.6882824681144086
$ws.Range("F12").Value = 0.7356546913071611
$ws.Range("G12").Value = 0.6483731122692404
$ws.Range("H12").Value = 0.7226104611102926
$ws.Range("I12").Value = 0.5630310507164111
$ws.Range("K12").Value = 0.6310311338324652
$ws.Range("B13").Value = 0.5521695235019877
$ws.Range("C13").Value = 0.1117098349841115
$ws.Range("E13").Value = 0.6846257058405172
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("G13").Value = 0.6480927613991128
$ws.Range("H13").Value = 0.722765950094626
$ws.Range("I13").Value = 0.5632143118593831
$ws.Range("K13").Value = 0.6278103233295269
$ws.Range("B14").Value = 0.5432042634148218
$ws.Range("C14").Value = 0.109784984603408
$ws.Range("E14").Value = 0.6727049047819094
$ws.Range("F14").Value = 0.7191683204515869
$ws.Range("G14").Value = 0.6471943648351868
$ws.Range("H14").Value = 0.7232847454952633
$ws.Range("I14").Value = 0.5638231203727528
$ws.Range("K14").Value = 0.6173069665047706
$ws.Range("B15").Value = 0.5377132555698267
$ws.Range("C15").Value = 0.1086058756370392
$ws.Range("E15").Value = 0.6654066402568048
$ws.Range("F15").Value = 0.7114413442032514
$ws.Range("G15").Value = 0.6466562629883867
$ws.Range("H15").Value = 0.7236115162376819
$ws.Range("I15").Value = 0.5642045962460358
$ws.Range("K15").Value = 0.610873670398604
$ws.Range("B16").Value = 0.5062539197415958
$ws.Range("C16").Value = 0.1018475314391196
$ws.Range("E16").Value = 0.6236364052578125
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("G16").Value = 0.6437587855223796
$ws.Range("H16").Value = 0.7256217010259718
$ws.Range("I16").Value = 0.5665216022705692
$ws.Range("K16").Value = 0.5740118977398367
$ws.Range("B17").Value = 0.4869614282214627
$ws.Range("C17").Value = 0.09770021489109126
$ws.Range("E17").Value = 0.5980574793645559
$ws.Range("F17").Value = 0.6400460337125793
$ws.Range("G17").Value = 0.6421466335094976
$ws.Range("H17").Value = 0.7269773870724805
$ws.Range("I17").Value = 0.5680594799843348
$ws.Range("K17").Value = 0.5514027712237919
$ws.Range("B18").Value = 0.4758674656762594
$ws.Range("C18").Value = 0.09531427857902486
$ws.Range("E18").Value = 0.5833613105617133
$ws.Range("F18").Value = 0.6244449056556647
$ws.Range("G18").Value = 0.6412798321435815
$ws.Range("H18").Value = 0.7278020525784541
$ws.Range("I18").Value = 0.5689866884105044
$ws.Range("K18").Value = 0.5384002301916837
$ws.Range("B19").Value = 0.4721116886341008
$ws.Range("C19").Value = 0.09450635073397962
$ws.Range("E19").Value = 0.5783881518655392
$ws.Range("F19").Value = 0.619163680173358
$ws.Range("G19").Value = 0.6409967045594698
$ws.Range("H19").Value = 0.7280889733526124
$ws.Range("I19").Value = 0.569307939186757
$ws.Range("K19").Value = 0.5339980754563669
$ws.Range("B20").Value = 0.4890148827151393
$ws.Range("C20").Value = 0.09814175576957496
$ws.Range("E20").Value = 0.600778713036064
$ws.Range("F20").Value = 0.642933953830422
$ws.Range("G20").Value = 0.6423119850470584
$ws.Range("H20").Value = 0.7268284214333249
$ws.Range("I20").Value = 0.567891351770804
$ws.Range("K20").Value = 0.5538093831003152
$ws.Range("B21").Value = 0.5458374794744714
$ws.Range("C21").Value = 0.1103503768188432
$ws.Range("E21").Value = 0.676205585233177
$ws.Range("F21").Value = 0.7228739723491628
$ws.Range("G21").Value = 0.6474557031974655
$ws.Range("H21").Value = 0.7231304878483797
$ws.Range("I21").Value = 0.5636425132328284
$ws.Range("K21").Value = 0.6203919918384599
$ws.Range("B22").Value = 0.5829956865332235
$ws.Range("C22").Value = 0.118325628069357
$ws.Range("E22").Value = 0.7256594144817967
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("G22").Value = 0.6513622913150101
$ws.Range("H22").Value = 0.7211159131950495
$ws.Range("I22").Value = 0.561248536940667
$ws.Range("K22").Value = 0.6639216521707851
$ws.Range("B23").Value = 0.563161769815764
$ws.Range("C23").Value = 0.1140693962934449
$ws.Range("E23").Value = 0.6992498261539168
$ws.Range("F23").Value = 0.7472568307830727
$ws.Range("G23").Value = 0.6492271038103752
$ws.Range("H23").Value = 0.722154212188812
$ws.Range("I23").Value = 0.5624910680794812
$ws.Range("K23").Value = 0.6406877847053067
$ws.Range("B24").Value = 0.4880865235629699
$ws.Range("C24").Value = 0.09794214009954771
$ws.Range("E24").Value = 0.5995484139147607
$ws.Range("F24").Value = 0.6416283278902171
$ws.Range("G24").Value = 0.6422370426380439
$ws.Range("H24").Value = 0.7268956278189336
$ws.Range("I24").Value = 0.5679672284240382
$ws.Range("K24").Value = 0.5527213671194602
$ws.Range("B25").Value = 0.4072045808032101
$ws.Range("C25").Value = 0.0805263154156961
$ws.Range("E25").Value = 0.4926132651867761
$ws.Range("F25").Value = 0.5279251897347166
$ws.Range("G25").Value = 0.6370050121206248
$ws.Range("H25").Value = 0.7337256793014149
$ws.Range("I25").Value = 0.5755034410775899
$ws.Range("K25").Value = 0.4578976842566931

Write-Host "Updated pl_mw result values for case with 380 kV"
